# Fixed Some Params Issues
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the DATA_TYPE value for param Q0019 (row 2) from "string" to "[0]string"
$ws.Range("E2").Value = "[0]string"

# Make sure the final selection / active cell is on E2, matching the authored edit
$ws.Activate()
$ws.Range("E2").Select()
